{"js": "// Locate the last bullet of the \"Semaine 6\" section: \"Faire un graphe avec\n// abscisse : instance id, ordonn\u00e9e : list des status.\" and append two new\n// bullet items (\"Pas facile ..\" and \"Test : Kibana 4\") right after it, each\n// keeping the same \"Paragraphedeliste\" list style / numbering as the\n// surrounding bullets.\nconst anchorText =\n  \"Faire un graphe avec abscisse : instance id, ordonn\u00e9e : list des status.\";\n\nconst results = context.document.body.search(anchorText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Anchor paragraph not found: \" + anchorText);\n}\n\nconst anchorParagraph = results.items[0].paragraphs.getFirst();\n\n// Inserting a literal newline right after the anchor's text splits the\n// paragraph, so the newly created paragraph(s) naturally inherit the same\n// paragraph style / bullet numbering (pStyle \"Paragraphedeliste\" +\n// numPr ilvl=0 / numId=1) as the paragraph being split, exactly like typing\n// at the end of the line and pressing Enter in Word.\nconst endOfAnchor = anchorParagraph.getRange(Word.RangeLocation.end);\nendOfAnchor.insertText(\"\\nPas facile ..\\nTest : Kibana 4\", Word.InsertLocation.end);\n\nawait context.sync();\n", "ps1": "# Locate the last bullet of the \"Semaine 6\" section: \"Faire un graphe avec\n# abscisse : instance id, ordonn\u00e9e : list des status.\" and append two new\n# bullet items (\"Pas facile ..\" and \"Test : Kibana 4\") right after it, each\n# keeping the same \"Paragraphedeliste\" list style / numbering (numId 1,\n# ilvl 0) as the surrounding bullets.\n\n$d = $word.ActiveDocument\n\n$anchorText = \"Faire un graphe avec abscisse : instance id, ordonn\u00e9e : list des status.\"\n\n$searchRange = $d.Content\n$found = $searchRange.Find.Execute($anchorText)\nif (-not $found) {\n    throw \"Anchor paragraph not found: $anchorText\"\n}\n\n$anchorParagraph = $searchRange.Paragraphs(1)\n\n# Build a zero-length range positioned right before the anchor paragraph's\n# end-of-paragraph mark (Range.End - 1), then insert the two new lines right\n# after that point. Because the insertion happens inside the existing list\n# paragraph (before its paragraph mark), Word splits it into three\n# paragraphs that all inherit the same paragraph style and list numbering\n# (pStyle \"Paragraphedeliste\", numPr ilvl=0/numId=1) as the original bullet.\n$insertionPoint = $d.Range($anchorParagraph.Range.End - 1, $anchorParagraph.Range.End - 1)\n$insertionPoint.InsertAfter(\"`rPas facile ..`rTest : Kibana 4\")\n"}
